$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 17 (shifting the existing rows 17-45 down to 20-48),
# picking up formatting (date style on column D) from the row above.
$ws.Range("A17:A19").EntireRow.Insert()

# New weekly data rows (Fecha serial 44547 = 2021-12-17)
$rows = @(
    @{ Row=17; Calidad="Extra";   Volumen=4000; PMin=1400; PMax=1500; PProm=1450 },
    @{ Row=18; Calidad="Primera"; Volumen=4000; PMin=950;  PMax=1000; PProm=975  },
    @{ Row=19; Calidad="Segunda"; Volumen=4000; PMin=800;  PMax=900;  PProm=850  }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Cells.Item($i, 1).Value = 8
    $ws.Cells.Item($i, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($i, 3).Value = "Coquimbo"
    $ws.Cells.Item($i, 4).Value = 44547
    $ws.Cells.Item($i, 5).Value = 4
    $ws.Cells.Item($i, 6).Value = 100112027
    $ws.Cells.Item($i, 7).Value = "Melón"
    $ws.Cells.Item($i, 8).Value = "Tuna"
    $ws.Cells.Item($i, 9).Value = $r.Calidad
    $ws.Cells.Item($i, 10).Value = $r.Volumen
    $ws.Cells.Item($i, 11).Value = $r.PMin
    $ws.Cells.Item($i, 12).Value = $r.PMax
    $ws.Cells.Item($i, 13).Value = $r.PProm
    $ws.Cells.Item($i, 14).Value = "$/unidad"
    $ws.Cells.Item($i, 15).Value = "Región de O'Higgins"
    $ws.Cells.Item($i, 16).Value = $r.PProm
    $ws.Cells.Item($i, 17).Value = 1
    $ws.Cells.Item($i, 18).Value = "Hortaliza"
}
